$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nine")

$ws.Range("E2").Value = 10.91

$ws.Range("D3").Value = 10.6
$ws.Range("E3").Value = 10.36

$ws.Range("C4").Value = 9.4
$ws.Range("E4").Value = 10.25
$ws.Range("G4").Value = 9.869999999999999

$ws.Range("B5").Value = 9.09
$ws.Range("C5").Value = 9.640000000000001
$ws.Range("D5").Value = 9.75
$ws.Range("F5").Value = 10.26

$ws.Range("E6").Value = 9.74
$ws.Range("G6").Value = 10.38
$ws.Range("H6").Value = 11.11
$ws.Range("J6").Value = 7.85

$ws.Range("D7").Value = 10.13
$ws.Range("F7").Value = 9.619999999999999
$ws.Range("I7").Value = 5.75

$ws.Range("F8").Value = 8.890000000000001

$ws.Range("G9").Value = 14.25

$ws.Range("F10").Value = 12.15
